$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.079.28"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").Value = "1.603.20"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3779"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3656"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.270"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.41%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08137"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.613"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001258"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.391"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.93%  "
$ws.Range("D17").Value = "1.600.75"
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06878"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.587"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.21%  "
$ws.Range("D24").Value = "23.088.24"
$ws.Range("E24").Value = "  -3.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.349"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.787"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.90%  "
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.287"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.373"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.827"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.92%  "
$ws.Range("D33").Value = "1.775.81"
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9558"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07685"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02724"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2547"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08906"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.369"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7094"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.06%  "
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07944"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.238"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.96%  "
